$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 79, pushing existing rows 79:112 down to 80:113.
$ws.Rows("79").Insert()

# Populate the new row 79 with its data (same dimension/category fields as the
# surrounding rows, new date + price/volume figures).
$ws.Range("A79").Value = 4
$ws.Range("B79").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C79").Value = "Los Lagos"
$ws.Range("D79").Value = 44825
$ws.Range("E79").Value = 10
$ws.Range("F79").Value = "Fruta"
$ws.Range("G79").Value = 100104
$ws.Range("H79").Value = "Frutos de pepita"
$ws.Range("I79").Value = 100104003
$ws.Range("J79").Value = "Membrillo"
$ws.Range("K79").Value = "Champion"
$ws.Range("L79").Value = "Primera"
$ws.Range("M79").Value = 40
$ws.Range("N79").Value = 14000
$ws.Range("O79").Value = 15000
$ws.Range("P79").Value = 14500
$ws.Range("Q79").Value = "`$/caja 18 kilos granel"
$ws.Range("R79").Value = "Región de O'Higgins"
$ws.Range("S79").Value = 806
$ws.Range("T79").Value = 18
